$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 361.89474
$ws.Range("I5").Value = 500.46155
$ws.Range("J5").Value = 61.666668
$ws.Range("K5").Value = 500.46155
$ws.Range("L5").Value = 61.666668
$ws.Range("M5").Value = -385.46155
$ws.Range("N5").Value = -291.666668
# Row 40
$ws.Range("H40").Value = 5185.4443
$ws.Range("I40").Value = 3139.8
$ws.Range("J40").Value = 7742.5
$ws.Range("K40").Value = 3139.8
$ws.Range("L40").Value = 7742.5
$ws.Range("M40").Value = -2964.8
$ws.Range("N40").Value = -8092.5
# Row 64
$ws.Range("H64").Value = 8421.714
$ws.Range("J64").Value = 8446.363
$ws.Range("L64").Value = 8446.363
$ws.Range("N64").Value = -8942.363
# Row 67
$ws.Range("H67").Value = 8421.714
$ws.Range("J67").Value = 8446.363
$ws.Range("L67").Value = 8446.363
$ws.Range("N67").Value = -10162.363
# Row 74
$ws.Range("H74").Value = 7432.952
$ws.Range("I74").Value = 5592.0
$ws.Range("J74").Value = 11114.857
$ws.Range("K74").Value = 5592.0
$ws.Range("L74").Value = 11114.857
$ws.Range("M74").Value = -4656.0
$ws.Range("N74").Value = -12986.857
# Row 77
$ws.Range("H77").Value = 7432.952
$ws.Range("I77").Value = 5592.0
$ws.Range("J77").Value = 11114.857
$ws.Range("K77").Value = 27960.0
$ws.Range("L77").Value = 55574.285
$ws.Range("M77").Value = -23280.0
$ws.Range("N77").Value = -64934.285
# Row 100
$ws.Range("H100").Value = 97306.63
$ws.Range("I100").Value = 148285.28
$ws.Range("J100").Value = 8094.0
$ws.Range("K100").Value = 148285.28
$ws.Range("L100").Value = 8094.0
$ws.Range("M100").Value = -147744.28
$ws.Range("N100").Value = -9176.0
# Row 113
$ws.Range("H113").Value = 61576.668
$ws.Range("I113").Value = 251624.5
$ws.Range("K113").Value = 251624.5
$ws.Range("M113").Value = -248370.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 119
$ws.Range("H119").Value = 63014.8
$ws.Range("J119").Value = 63014.8
$ws.Range("L119").Value = 63014.8
$ws.Range("N119").Value = -72690.8
# Row 125
$ws.Range("H125").Value = 48306.145
$ws.Range("J125").Value = 48306.145
$ws.Range("L125").Value = 48306.145
$ws.Range("N125").Value = -58146.145
# Row 135
$ws.Range("H135").Value = 47264.11
$ws.Range("J135").Value = 47264.11
$ws.Range("L135").Value = 47264.11
$ws.Range("N135").Value = -57404.11
# Row 137
$ws.Range("H137").Value = 69991.0
$ws.Range("J137").Value = 69991.0
$ws.Range("L137").Value = 69991.0
$ws.Range("N137").Value = -80191.0
# Row 139
$ws.Range("H139").Value = 60542.0
$ws.Range("J139").Value = 60542.0
$ws.Range("L139").Value = 60542.0
$ws.Range("N139").Value = -70822.0
# Row 140
$ws.Range("H140").Value = 69130.836
$ws.Range("I140").Value = 0.0
$ws.Range("J140").Value = 69130.836
$ws.Range("K140").Value = 0.0
$ws.Range("L140").Value = 69130.836
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -79490.836
# Row 141
$ws.Range("H141").Value = 49800.0
$ws.Range("J141").Value = 49800.0
$ws.Range("L141").Value = 49800.0
$ws.Range("N141").Value = -60160.0

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 768.0
$ws.Range("I22").Value = 768.0
$ws.Range("J22").Value = 0.0
$ws.Range("K22").Value = 768.0
$ws.Range("L22").Value = 0.0
$ws.Range("M22").Value = -595.0
$ws.Range("N22").ClearContents()
# Row 132
$ws.Range("H132").Value = 65995.0
$ws.Range("J132").Value = 65995.0
$ws.Range("L132").Value = 65995.0
$ws.Range("N132").Value = -76115.0
# Row 133
$ws.Range("H133").Value = 0.0
$ws.Range("J133").Value = 0.0
$ws.Range("L133").Value = 0.0
$ws.Range("N133").ClearContents()
# Row 134
$ws.Range("H134").Value = 3102.5
$ws.Range("I134").Value = 1779.0
$ws.Range("K134").Value = 5337.0
$ws.Range("M134").Value = -2802.0
# Row 135
$ws.Range("H135").Value = 45410.0
$ws.Range("J135").Value = 45410.0
$ws.Range("L135").Value = 45410.0
$ws.Range("N135").Value = -55550.0
# Row 137
$ws.Range("H137").Value = 69991.4
$ws.Range("J137").Value = 69991.4
$ws.Range("L137").Value = 69991.4
$ws.Range("N137").Value = -80191.4
# Row 138
$ws.Range("H138").Value = 66996.0
$ws.Range("J138").Value = 66996.0
$ws.Range("L138").Value = 66996.0
$ws.Range("N138").Value = -77276.0

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 6138.231
$ws.Range("I58").Value = 2224.7144
$ws.Range("J58").Value = 10704.0
$ws.Range("K58").Value = 2224.7144
$ws.Range("L58").Value = 10704.0
$ws.Range("M58").Value = -2021.7144
$ws.Range("N58").Value = -11110.0
# Row 99
$ws.Range("H99").Value = 3850.0
$ws.Range("I99").Value = 3850.0
$ws.Range("K99").Value = 3850.0
$ws.Range("M99").Value = -2352.0
# Row 126
$ws.Range("H126").Value = 3850.0
$ws.Range("I126").Value = 3850.0
$ws.Range("K126").Value = 11550.0
$ws.Range("M126").Value = -9080.0
# Row 132
$ws.Range("H132").Value = 10671.333
$ws.Range("I132").Value = 12000.0
$ws.Range("J132").Value = 10007.0
$ws.Range("K132").Value = 36000.0
$ws.Range("L132").Value = 30021.0
$ws.Range("M132").Value = -33470.0
$ws.Range("N132").Value = -35081.0
# Row 133
$ws.Range("H133").Value = 56497.2
$ws.Range("I133").Value = 21999.0
$ws.Range("J133").Value = 60330.332
$ws.Range("K133").Value = 21999.0
$ws.Range("L133").Value = 60330.332
$ws.Range("M133").Value = -19469.0
$ws.Range("N133").Value = -65390.332
# Row 134
$ws.Range("H134").Value = 8974.75
$ws.Range("I134").Value = 2000.0
$ws.Range("J134").Value = 11299.667
$ws.Range("K134").Value = 6000.0
$ws.Range("L134").Value = 33899.001
$ws.Range("M134").Value = -3465.0
$ws.Range("N134").Value = -38969.001
# Row 135
$ws.Range("H135").Value = 69999.0
$ws.Range("J135").Value = 69999.0
$ws.Range("L135").Value = 69999.0
$ws.Range("N135").Value = -80139.0
# Row 136
$ws.Range("H136").Value = 6138.231
$ws.Range("I136").Value = 2224.7144
$ws.Range("J136").Value = 10704.0
$ws.Range("K136").Value = 6674.1432
$ws.Range("L136").Value = 32112.0
$ws.Range("M136").Value = -4124.1432
$ws.Range("N136").Value = -37212.0
# Row 137
$ws.Range("H137").Value = 79119.8
$ws.Range("J137").Value = 79119.8
$ws.Range("L137").Value = 79119.8
$ws.Range("N137").Value = -89319.8
# Row 138
$ws.Range("H138").Value = 68499.5
$ws.Range("J138").Value = 68499.5
$ws.Range("L138").Value = 68499.5
$ws.Range("N138").Value = -78779.5
# Row 139
$ws.Range("H139").Value = 58000.0
$ws.Range("J139").Value = 58000.0
$ws.Range("L139").Value = 58000.0
$ws.Range("N139").Value = -68280.0
# Row 140
$ws.Range("H140").Value = 64880.332
$ws.Range("J140").Value = 64880.332
$ws.Range("L140").Value = 64880.332
$ws.Range("N140").Value = -75240.332
# Row 141
$ws.Range("H141").Value = 233994.8
$ws.Range("J141").Value = 233994.8
$ws.Range("L141").Value = 233994.8
$ws.Range("N141").Value = -244354.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 6477681.0
$ws.Range("I4").Value = 4210609.0
$ws.Range("K4").Value = 12631827.0
$ws.Range("M4").Value = -12631715.0
# Row 5
$ws.Range("H5").Value = 16921.875
$ws.Range("I5").Value = 2063.5
$ws.Range("J5").Value = 31780.25
$ws.Range("K5").Value = 6190.5
$ws.Range("L5").Value = 95340.75
$ws.Range("M5").Value = -6078.5
$ws.Range("N5").Value = -95564.75
# Row 50
$ws.Range("H50").Value = 37045520.0
$ws.Range("I50").Value = 55557380.0
$ws.Range("J50").Value = 21800.0
$ws.Range("K50").Value = 166672140.0
$ws.Range("L50").Value = 65400.0
$ws.Range("M50").Value = -166671659.0
$ws.Range("N50").Value = -66362.0
# Row 53
$ws.Range("H53").Value = 37045520.0
$ws.Range("I53").Value = 55557380.0
$ws.Range("J53").Value = 21800.0
$ws.Range("K53").Value = 166672140.0
$ws.Range("L53").Value = 65400.0
$ws.Range("M53").Value = -166671659.0
$ws.Range("N53").Value = -66362.0
# Row 122
$ws.Range("H122").Value = 6482450.0
$ws.Range("J122").Value = 6804768.5
$ws.Range("L122").Value = 61242916.5
$ws.Range("N122").Value = -61247816.5
# Row 135
$ws.Range("H135").Value = 16921.875
$ws.Range("I135").Value = 2063.5
$ws.Range("J135").Value = 31780.25
$ws.Range("K135").Value = 18571.5
$ws.Range("L135").Value = 286022.25
$ws.Range("M135").Value = -16036.5
$ws.Range("N135").Value = -291092.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 49
$ws.Range("H49").Value = 23999.834
$ws.Range("I49").Value = 12500.0
$ws.Range("J49").Value = 29749.75
$ws.Range("K49").Value = 12500.0
$ws.Range("L49").Value = 29749.75
$ws.Range("M49").Value = -12316.0
$ws.Range("N49").Value = -30117.75
# Row 117
$ws.Range("H117").Value = 65000.0
$ws.Range("J117").Value = 65000.0
$ws.Range("L117").Value = 65000.0
$ws.Range("N117").Value = -71884.0
# Row 132
$ws.Range("H132").Value = 550006.4
$ws.Range("I132").Value = 627292.6
$ws.Range("K132").Value = 1881877.8
$ws.Range("M132").Value = -1879347.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 2489.0
$ws.Range("J46").Value = 2460.611
$ws.Range("L46").Value = 2460.611
$ws.Range("N46").Value = -2836.611
# Row 50
$ws.Range("H50").Value = 31500.0
$ws.Range("J50").Value = 31500.0
$ws.Range("L50").Value = 31500.0
$ws.Range("N50").Value = -32774.0
# Row 54
$ws.Range("H54").Value = 40000.0
$ws.Range("J54").Value = 40000.0
$ws.Range("L54").Value = 40000.0
$ws.Range("N54").Value = -41288.0
# Row 127
$ws.Range("H127").Value = 0.0
$ws.Range("J127").Value = 0.0
$ws.Range("L127").Value = 0.0
$ws.Range("N127").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Range("H70").Value = 32249.125
$ws.Range("I70").Value = 33399.4
$ws.Range("J70").Value = 30332.0
$ws.Range("K70").Value = 33399.4
$ws.Range("L70").Value = 30332.0
$ws.Range("M70").Value = -33084.4
$ws.Range("N70").Value = -30962.0
# Row 73
$ws.Range("H73").Value = 32249.125
$ws.Range("I73").Value = 33399.4
$ws.Range("J73").Value = 30332.0
$ws.Range("K73").Value = 33399.4
$ws.Range("L73").Value = 30332.0
$ws.Range("M73").Value = -32307.4
$ws.Range("N73").Value = -32516.0
# Row 122
$ws.Range("H122").Value = 3658.2778
$ws.Range("I122").Value = 1769.2667
$ws.Range("J122").Value = 13103.333
$ws.Range("K122").Value = 5307.800099999999
$ws.Range("L122").Value = 39309.999
$ws.Range("M122").Value = -2857.800099999999
$ws.Range("N122").Value = -44209.999
